# Cosmetic edit: update the Status value in D2 from PASS to FAIL
# and move the active selection to F10 (matches the saved cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "FAIL"

$ws.Range("F10").Select()
